$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Séparer le cours "Routes, contexte et API" en deux semaines distinctes :
# la semaine 16 (ligne 17) ne garde que les Routes, et une nouvelle semaine
# (ligne 18) couvre le Contexte/API + l'Internationalisation/Accessibilité,
# ce qui décale Authentification à l'exercice 13.

$ws.Range("D17").Value = "[Exercice 10 - Routes](exercice10_routes.md)"

$ws.Range("C18").Value = "[Contexte et API](react4.md)<br/> [Régles de React](regles_react.md)<br/> [Internationalisation](internationalisation.md) <br />[Accessibilité](accessibilite.md)"
$ws.Range("D18").Value = "[Exercice 11 - Contexte et API](exercice11_context.md)<br/>[Exercice 12 - Internatialisation](exercice12_internationalisation.md)"

$ws.Range("D19").Value = "[Exercice 13 - Authentification](exercice13_authentification.md)"

# L'auteur a fait défiler la feuille et sélectionné la nouvelle ligne ajoutée.
$ws.Range("D18").Select()
